# codelijst aanvullen met altLabel
#
# Inserts a new "altLabel" column (between "theme" and "broader") into the
# verwaarlozing codelijst worksheet, populates it, and removes the trailing
# period from the "prefLabel" values of the individual concept rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at H ("theme" is G, old "broader" was H).
#    This shifts the old H:R block to I:S and keeps all existing values intact.
$ws.Columns("H").Insert()

# 2. Header for the newly inserted column.
$ws.Range("H1").Value = "altLabel"

# 3. Populate the altLabel column per row.
#    - Collection (row 2) and ConceptScheme (row 16) rows have no altLabel.
$ws.Range("H2").Value = "null"
$ws.Range("H16").Value = "null"

#    - Leaf concepts (with a "broader" concept) get the capitalized notation
#      as their altLabel.
$ws.Range("H3").Value = "Buitenmuren"
$ws.Range("H4").Value = "Buitentimmerwerk"
$ws.Range("H5").Value = "Dakbedekking"
$ws.Range("H6").Value = "Dakgebinte"
$ws.Range("H7").Value = "Dakgoten"
$ws.Range("H9").Value = "Kroonlijst"
$ws.Range("H10").Value = "Liften"
$ws.Range("H11").Value = "Schoorstenen"
$ws.Range("H13").Value = "Trappen"

#    - Top concepts without a "broader" (rows 8, 12, 14, 15) reuse their own
#      (period-stripped) prefLabel as the altLabel.
$ws.Range("H8").Value = "Gebrek aan de toestand van een gebouwonderdeel"
$ws.Range("H12").Value = "Gebrek dat stabiliteit in het gedrang brengt"
$ws.Range("H14").Value = "Gebrek dat veiligheid in het gedrang brengt"
$ws.Range("H15").Value = "Gebrek dat leidt tot vochtindringing"

# 4. Strip the trailing period from prefLabel (column F) for the concept rows.
$periodRows = 3,4,5,6,7,8,9,10,11,12,13,14,15
foreach ($r in $periodRows) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val.EndsWith(".")) {
        $cell.Value = $val.TrimEnd(".")
    }
}
